# Applies:
#  1) Slide 4 notes: prepend "Sources: ..." citation block before the existing speaker note.
#  2) Slide 5 notes: prepend "Source: Addy Osmani - URL" line before the existing speaker note.
#  3) Slide 4 body ("DataReportal 2026, Reuters" caption): split into three runs so that
#     "DataReportal 2026" and "Reuters" become clickable, underlined hyperlinks, keeping the
#     ", " separator plain.
#  4) Slide 5 body ("Addy Osmani" caption): turn the whole caption into a clickable, underlined
#     hyperlink.

$p = $ppt.ActivePresentation

$nl = [char]10
$bullet = [char]0x2022

# ---------------------------------------------------------------------------
# 1) Slide 4 speaker notes - add source citations above the existing note text
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$notes4 = $slide4.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notes4Line1 = "Sources: DataReportal 2026, Reuters."
$notes4Line2 = $bullet + " https://datareportal.com/reports/digital-2026-one-billion-people-using-ai"
$notes4Line3 = $bullet + " https://www.reuters.com/commentary/breakingviews/ai-investment-bubble-inflated-by-trio-dilemmas-2025-09-25/"
$notes4Line4 = "Mass adoption but shallow depth. The opportunity gap is enormous. ~15 seconds."
$notes4.Text = $notes4Line1 + $nl + $notes4Line2 + $nl + $notes4Line3 + $nl + $notes4Line4

# ---------------------------------------------------------------------------
# 2) Slide 5 speaker notes - add source citation above the existing note text
# ---------------------------------------------------------------------------
$emdash = [char]0x2014
$rsquo = [char]0x2019
$slide5 = $p.Slides.Item(5)
$notes5 = $slide5.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notes5Line1 = "Source: Addy Osmani " + $emdash + " https://addyo.substack.com/p/the-reality-of-ai-assisted-software"
$notes5Line2 = "Developers who felt 20% faster actually took 19% longer once debugging was included. The gap isn" + $rsquo + "t just unused potential " + $emdash + " it" + $rsquo + "s active harm. ~15 seconds."
$notes5.Text = $notes5Line1 + $nl + $notes5Line2

# ---------------------------------------------------------------------------
# 3) Slide 4 caption "DataReportal 2026, Reuters" -> hyperlinked source names
# ---------------------------------------------------------------------------
$caption4 = $slide4.Shapes.Item(4).TextFrame.TextRange
# "DataReportal 2026" = chars 1-17
$dataReportalRun = $caption4.Characters(1, 17)
$dataReportalRun.Font.Underline = 1
$dataReportalRun.ActionSettings(1).Hyperlink.Address = "https://datareportal.com/reports/digital-2026-one-billion-people-using-ai"
# "Reuters" = chars 20-26
$reutersRun = $caption4.Characters(20, 7)
$reutersRun.Font.Underline = 1
$reutersRun.ActionSettings(1).Hyperlink.Address = "https://www.reuters.com/commentary/breakingviews/ai-investment-bubble-inflated-by-trio-dilemmas-2025-09-25/"

# ---------------------------------------------------------------------------
# 4) Slide 5 caption "Addy Osmani" -> hyperlinked source name
# ---------------------------------------------------------------------------
$caption5 = $slide5.Shapes.Item(4).TextFrame.TextRange
$caption5.Font.Underline = 1
$caption5.ActionSettings(1).Hyperlink.Address = "https://addyo.substack.com/p/the-reality-of-ai-assisted-software"
